# Updated cryptos list on Tue Jul 16 20:29:17 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.961.87'
$ws.Range("E2").Value = '  +2.20%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.463.26'
$ws.Range("E3").Value = '  +0.78%  '
# Row 4
$ws.Range("E4").Value = '  +0.05%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.61'
$ws.Range("E5").Value = '  -0.21%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.65'
$ws.Range("E6").Value = '  +2.56%  '
# Row 7
$ws.Range("E7").Value = '  +0.12%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.464.38'
$ws.Range("E8").Value = '  +0.72%  '
# Row 9
$ws.Range("E9").Value = '  +7.75%  '
# Row 10
$ws.Range("E10").Value = '  -2.51%  '
# Row 11
$ws.Range("E11").Value = '  +2.45%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.439'
$ws.Range("E12").Value = '  +0.66%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.062.26'
$ws.Range("E13").Value = '  +1.09%  '
# Row 14
$ws.Range("E14").Value = '  -2.46%  '
# Row 15
$ws.Range("E15").Value = '  +4.48%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.17'
$ws.Range("E16").Value = '  +3.50%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.992.67'
$ws.Range("E17").Value = '  +1.96%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.501.15'
$ws.Range("E18").Value = '  +2.73%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.37'
$ws.Range("E19").Value = '  -1.07%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.31'
$ws.Range("E20").Value = '  +0.46%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '388.68'
$ws.Range("E21").Value = '  -0.94%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.22'
$ws.Range("E22").Value = '  -3.33%  '
# Row 23
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.20'
$ws.Range("E23").Value = '  +1.91%  '
# Row 24
$ws.Range("B24").Value = 'Polygon'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.545'
$ws.Range("E24").Value = '  +0.94%  '
# Row 25
$ws.Range("E25").Value = '  +0.13%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000124'
$ws.Range("E26").Value = '  +15.28%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.73'
$ws.Range("E27").Value = '  +1.83%  '
# Row 28
$ws.Range("E28").Value = '  -0.09%  '
# Row 29
$ws.Range("E29").Value = '  +0.05%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.24'
$ws.Range("E30").Value = '  +8.72%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.43'
$ws.Range("E31").Value = '  +6.06%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.04'
$ws.Range("E32").Value = '  +0.10%  '
# Row 33
$ws.Range("E33").Value = '  +0.87%  '
# Row 34
$ws.Range("E34").Value = '  -0.78%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  +0.07%  '
# Row 36
$ws.Range("E36").Value = '  +4.59%  '
# Row 37
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.24'
$ws.Range("E37").Value = '  +2.89%  '
# Row 38
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.51'
$ws.Range("E38").Value = '  +0.66%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.93'
$ws.Range("E39").Value = '  +2.48%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.000.87'
$ws.Range("E40").Value = '  +2.16%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0768'
$ws.Range("E41").Value = '  -1.58%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.29'
$ws.Range("E42").Value = '  -2.66%  '
# Row 43
$ws.Range("E43").Value = '  +4.94%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '43.00'
$ws.Range("E44").Value = '  +3.13%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0316'
$ws.Range("E45").Value = '  -1.87%  '
# Row 46
$ws.Range("E46").Value = '  +1.03%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.42'
$ws.Range("E47").Value = '  +8.36%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.09'
$ws.Range("E48").Value = '  +1.00%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.878'
$ws.Range("E49").Value = '  +7.40%  '
# Row 50
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.58'
$ws.Range("E50").Value = '  +3.24%  '
# Row 51
$ws.Range("B51").Value = 'Bittensor'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '305.36'
$ws.Range("E51").Value = '  +2.85%  '